$d = $word.ActiveDocument

$replacements = @{
    "144×9=" = "900×2=";
    "114×2=" = "499×7=";
    "726×7=" = "354×5=";
    "261×4=" = "564×8=";
    "304×4=" = "177×9=";
    "523×9=" = "956×6=";
    "572×4=" = "890×9=";
    "749×9=" = "373×8=";
    "103×7=" = "331×6=";
    "719×9=" = "922×3=";
    "655×4=" = "292×2=";
    "521×5=" = "742×4=";
    "237×3=" = "526×9=";
    "124×2=" = "955×3=";
    "910×5=" = "799×6=";
    "918×8=" = "107×8=";
    "194×5=" = "477×9=";
    "303×4=" = "343×8=";
    "715×9=" = "958×9=";
    "367×4=" = "266×6=";
    "382×6=" = "112×6=";
    "927×5=" = "783×4=";
    "595×7=" = "589×3=";
    "722×8=" = "809×6=";
    "880×4=" = "718×5=";
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
